# Scheduled runner update: refresh market-board derived price/profit
# columns (H:N) on the Golem_Profits leve-crafting sheets with the
# latest pulled values. Columns: H=currentAveragePrice,
# I=currentAveragePriceNQ, J=currentAveragePriceHQ, K=LevePriceNQ,
# L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 742
$ws.Cells.Item(38, 9).Value = 399.18182
$ws.Cells.Item(38, 10).Value = 1999
$ws.Cells.Item(38, 11).Value = 1197.54546
$ws.Cells.Item(38, 12).Value = 5997
$ws.Cells.Item(38, 13).Value = -825.54546
$ws.Cells.Item(38, 14).Value = -6741
$ws.Cells.Item(39, 8).Value = 135.28572
$ws.Cells.Item(39, 9).Value = 135.28572
$ws.Cells.Item(39, 11).Value = 405.85716
$ws.Cells.Item(39, 13).Value = -109.85716
$ws.Cells.Item(40, 8).Value = 1638.3077
$ws.Cells.Item(40, 9).Value = 1627.2727
$ws.Cells.Item(40, 10).Value = 1699
$ws.Cells.Item(40, 11).Value = 1627.2727
$ws.Cells.Item(40, 12).Value = 1699
$ws.Cells.Item(40, 13).Value = -1452.2727
$ws.Cells.Item(40, 14).Value = -2049
$ws.Cells.Item(41, 8).Value = 216
$ws.Cells.Item(41, 9).Value = 180.8
$ws.Cells.Item(41, 10).Value = 274.66666
$ws.Cells.Item(41, 11).Value = 180.8
$ws.Cells.Item(41, 12).Value = 274.66666
$ws.Cells.Item(41, 13).Value = 259.2
$ws.Cells.Item(41, 14).Value = -1154.66666
$ws.Cells.Item(58, 8).Value = 3013.75
$ws.Cells.Item(58, 9).Value = 1703.3334
$ws.Cells.Item(58, 10).Value = 3800
$ws.Cells.Item(58, 11).Value = 5110.0002
$ws.Cells.Item(58, 12).Value = 11400
$ws.Cells.Item(58, 13).Value = -4960.0002
$ws.Cells.Item(58, 14).Value = -11700
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 999.5
$ws.Cells.Item(2, 9).Value = 999.5
$ws.Cells.Item(2, 11).Value = 999.5
$ws.Cells.Item(2, 13).Value = -886.5
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 999.5
$ws.Cells.Item(116, 9).Value = 999.5
$ws.Cells.Item(116, 11).Value = 999.5
$ws.Cells.Item(116, 13).Value = 1294.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 999.5
$ws.Cells.Item(3, 9).Value = 999.5
$ws.Cells.Item(3, 11).Value = 999.5
$ws.Cells.Item(3, 13).Value = -885.5
$ws.Cells.Item(20, 8).Value = 973.5
$ws.Cells.Item(20, 9).Value = 965.6667
$ws.Cells.Item(20, 10).Value = 997
$ws.Cells.Item(20, 11).Value = 965.6667
$ws.Cells.Item(20, 12).Value = 997
$ws.Cells.Item(20, 13).Value = -718.6667
$ws.Cells.Item(20, 14).Value = -1491
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(88, 8).Value = 19791.857
$ws.Cells.Item(88, 10).Value = 20923.834
$ws.Cells.Item(88, 12).Value = 20923.834
$ws.Cells.Item(88, 14).Value = -21735.834
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(91, 8).Value = 19791.857
$ws.Cells.Item(91, 10).Value = 20923.834
$ws.Cells.Item(91, 12).Value = 20923.834
$ws.Cells.Item(91, 14).Value = -23731.834
$ws.Cells.Item(106, 8).Value = 2766.6667
$ws.Cells.Item(106, 10).Value = 2766.6667
$ws.Cells.Item(106, 12).Value = 2766.6667
$ws.Cells.Item(106, 14).Value = -5290.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 25000
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(62, 8).Value = 4500
$ws.Cells.Item(62, 9).Value = 4500
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 4500
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -3876
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 4500
$ws.Cells.Item(65, 9).Value = 4500
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 22500
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -19380
$ws.Cells.Item(65, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 383.66666
$ws.Cells.Item(7, 9).Value = 401
$ws.Cells.Item(7, 10).Value = 375
$ws.Cells.Item(7, 11).Value = 1203
$ws.Cells.Item(7, 12).Value = 1125
$ws.Cells.Item(7, 13).Value = -1091
$ws.Cells.Item(7, 14).Value = -1349
$ws.Cells.Item(34, 8).Value = 4009.3635
$ws.Cells.Item(34, 9).Value = 1600
$ws.Cells.Item(34, 10).Value = 4544.778
$ws.Cells.Item(34, 11).Value = 4800
$ws.Cells.Item(34, 12).Value = 13634.334
$ws.Cells.Item(34, 13).Value = -4716
$ws.Cells.Item(34, 14).Value = -13802.334
$ws.Cells.Item(39, 8).Value = 13645.111
$ws.Cells.Item(39, 10).Value = 13645.111
$ws.Cells.Item(39, 12).Value = 40935.333
$ws.Cells.Item(39, 14).Value = -41523.333
$ws.Cells.Item(51, 8).Value = 1409.6666
$ws.Cells.Item(51, 9).Value = 1369.5
$ws.Cells.Item(51, 11).Value = 4108.5
$ws.Cells.Item(51, 13).Value = -3648.5
$ws.Cells.Item(55, 8).Value = 2526.5557
$ws.Cells.Item(55, 10).Value = 2879.8518
$ws.Cells.Item(55, 12).Value = 8639.555399999999
$ws.Cells.Item(55, 14).Value = -8993.555399999999
$ws.Cells.Item(121, 8).Value = 419.16666
$ws.Cells.Item(121, 9).Value = 419.16666
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 1257.49998
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = 52.50001999999995
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 4460.5557
$ws.Cells.Item(131, 10).Value = 4810.8335
$ws.Cells.Item(131, 12).Value = 14432.5005
$ws.Cells.Item(131, 14).Value = -24512.5005
$ws.Cells.Item(134, 8).Value = 2097.5
$ws.Cells.Item(134, 9).Value = 2097.5
$ws.Cells.Item(134, 11).Value = 6292.5
$ws.Cells.Item(134, 13).Value = -1222.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6466.3335
$ws.Cells.Item(70, 10).Value = 5949.5
$ws.Cells.Item(70, 12).Value = 5949.5
$ws.Cells.Item(70, 14).Value = -6489.5
$ws.Cells.Item(73, 8).Value = 6466.3335
$ws.Cells.Item(73, 10).Value = 5949.5
$ws.Cells.Item(73, 12).Value = 5949.5
$ws.Cells.Item(73, 14).Value = -7821.5
$ws.Cells.Item(132, 8).Value = 909.6667
$ws.Cells.Item(132, 9).Value = 909.6667
$ws.Cells.Item(132, 11).Value = 2729.0001
$ws.Cells.Item(132, 13).Value = -199.0001000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 1731.125
$ws.Cells.Item(32, 10).Value = 2500
$ws.Cells.Item(32, 12).Value = 2500
$ws.Cells.Item(32, 14).Value = -3134
$ws.Cells.Item(46, 8).Value = 2558.6
$ws.Cells.Item(46, 9).Value = 698.25
$ws.Cells.Item(46, 11).Value = 698.25
$ws.Cells.Item(46, 13).Value = -510.25
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 3028.6667
$ws.Cells.Item(81, 9).Value = 3258.5
$ws.Cells.Item(81, 10).Value = 1190
$ws.Cells.Item(81, 11).Value = 6517
$ws.Cells.Item(81, 12).Value = 2380
$ws.Cells.Item(81, 13).Value = -5456
$ws.Cells.Item(81, 14).Value = -4502
$ws.Cells.Item(84, 8).Value = 3028.6667
$ws.Cells.Item(84, 9).Value = 3258.5
$ws.Cells.Item(84, 10).Value = 1190
$ws.Cells.Item(84, 11).Value = 32585
$ws.Cells.Item(84, 12).Value = 11900
$ws.Cells.Item(84, 13).Value = -27281
$ws.Cells.Item(84, 14).Value = -22508
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()
